$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new cell text, derived from the target diff of the
# underlying OOXML (coin price/volume refresh + an ARBITRUM/InjectiveProtocol
# row swap). All source cells are plain text (t="inlineStr"), including
# price cells that look numeric (e.g. "226.20"), so those are entered with a
# leading apostrophe to stop Excel from coercing them into numbers.
$updates = @{
    'D2' = '34.372.72'
    'E2' = '  +0.86%  '
    'D3' = '1.787.60'
    'E3' = '  +0.36%  '
    'E4' = '  -0.08%  '
    'D5' = '''226.20'
    'E5' = '  +0.36%  '
    'D6' = '''0.555'
    'E6' = '  +1.90%  '
    'E7' = '  -0.08%  '
    'D8' = '''32.95'
    'E8' = '  +3.51%  '
    'E9' = '  +1.20%  '
    'D11' = '''0.0946'
    'E11' = '  -0.10%  '
    'D12' = '2.046.25'
    'E12' = '  +0.43%  '
    'D13' = '''11.16'
    'E13' = '  +2.10%  '
    'D14' = '1.787.89'
    'E14' = '  +0.10%  '
    'D15' = '''0.634'
    'E15' = '  +2.09%  '
    'D16' = '34.318.89'
    'E16' = '  +0.65%  '
    'E17' = '  +2.65%  '
    'D18' = '''68.45'
    'E18' = '  +1.33%  '
    'D19' = '''245.12'
    'E19' = '  -0.22%  '
    'D20' = '0.0₃0794'
    'E20' = '  +0.87%  '
    'D21' = '''11.22'
    'E21' = '  +3.31%  '
    'E22' = '  -0.13%  '
    'E23' = '  +1.39%  '
    'D24' = '''168.49'
    'E24' = '  +4.23%  '
    'E25' = '  +1.44%  '
    'E26' = '  +3.13%  '
    'D27' = '''16.57'
    'E27' = '  +1.96%  '
    'E28' = '  +1.39%  '
    'E30' = '  +8.27%  '
    'E31' = '  +1.93%  '
    'E32' = '  +2.58%  '
    'E33' = '  +0.84%  '
    'D35' = '''2.61'
    'E35' = '  +6.33%  '
    'D36' = '1.409.28'
    'E36' = '  -2.87%  '
    'D37' = '''0.681'
    'E37' = '  +4.89%  '
    'E38' = '  +2.95%  '
    'E39' = '  +0.41%  '
    'D40' = '''84.30'
    'E40' = '  +4.57%  '
    'D41' = '''2.79'
    'E41' = '  +2.68%  '
    'E42' = '  +0.34%  '
    'B43' = 'ARBITRUM'
    'C43' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D43' = '''0.938'
    'E43' = '  +2.62%  '
    'B44' = 'InjectiveProtocol'
    'C44' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'D44' = '''13.98'
    'E44' = '  +1.37%  '
    'E45' = '  +2.00%  '
    'D46' = '''1.11'
    'E46' = '  +2.84%  '
    'E47' = '  +0.56%  '
    'D48' = '1.946.74'
    'E48' = '  +0.46%  '
    'D49' = '''105.41'
    'E49' = '  +0.97%  '
    'E50' = '  -0.09%  '
    'D51' = '0.0₆0129'
    'E51' = '  -3.29%  '
}

# Cells whose value needed the text-forcing apostrophe above; once the value
# is safely stored as text we drop back to the default "Normal" style so we
# don't leave a stray quote-prefixed style behind.
$textForced = @('D5', 'D6', 'D8', 'D11', 'D13', 'D15', 'D18', 'D19', 'D21', 'D24', 'D27', 'D35', 'D37', 'D40', 'D41', 'D43', 'D44', 'D46', 'D49')

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in $textForced) {
    $ws.Range($ref).Style = "Normal"
}
